# Security Feedback 3a: Add validation for email to be not empty
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("3. Security")
$ws.Activate()

# Insert a new row at row 7, shifting existing rows 7+ down by one
$ws.Rows("7:7").Insert()

# Populate the new row with the feedback text
$ws.Range("C7").Value = "Security Feedback 3a:"
$ws.Range("D7").Value = "Add validation for email to be not empty"

# Match the author's final selection/view state (cursor moved to D8 after insert)
$ws.Range("D8").Select()
